$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# More Redfin test case cities, appended after the existing list of cities.
$newCities = @("Portland", "Sacramento", "Baltimore", "Milwaukee", "Detroit", "Minneapolis")

# Find the last populated row in column A (xlUp = -4162) and append below it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$row = $lastRow + 1
foreach ($city in $newCities) {
    $ws.Range("A$row").Value = $city
    $row++
}

$lastNewRow = $row - 1
$ws.Range("A$lastNewRow").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
